$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""
